$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates for the crypto price/volume refresh.
# Force text number format on changed Price/Volume cells so values
# (which look numeric, e.g. "219.37") are kept as text like the original inline strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.251.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.662.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.37"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2669"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06344"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.09"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07725"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.693.96"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.436"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.891.68"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5481"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8239"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.274.97"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.660"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.77"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.15"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1246"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.240"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.20"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05975"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.90%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.611"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.296"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.59%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9837"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.425"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.775"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5898"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.027"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01599"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8599"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.030.02"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.13"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.805.24"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.47"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.015"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.072"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05185"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.468"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.37%  "
